# Generate Report for Archive
# Update localization status for the two files that moved from
# "Ready for handoff" to "In Translation": 4df4343e-... and 8b3ea37d-...
# This touches the Overview sheet (per-language status columns) as well
# as the per-language detail sheets (Status column).

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
# Row 3 -> 4df4343e-8da6-4563-b670-4855371533e4.md
$wsOverview.Range("E3").Value = "In Translation"
$wsOverview.Range("F3").Value = "In Translation"
# Row 4 -> 8b3ea37d-c86c-48ba-beda-d897d9c068ee.md
$wsOverview.Range("E4").Value = "In Translation"
$wsOverview.Range("F4").Value = "In Translation"

# --- zh-cn detail sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "In Translation"
$wsZhCn.Range("C4").Value = "In Translation"

# --- de-de detail sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "In Translation"
$wsDeDe.Range("C4").Value = "In Translation"
